# "did survey age compe for 2 sex models"
# Update the "Things to update for sex model" tracker sheet:
# - Clear the old E1/E2 note ("Bigger things" / "srv_comp_hat") that lived
#   in the header area, and write a fresh "Questions" / "sex ratio in the
#   likelihood?" note at E3:E4 instead.
# - Add a "y" Done-flag in column B for rows 10-18 (the survey/fishery
#   composition rows), and attach a note column (E) describing what was
#   done for each.
# - Add a new block of notes in columns E:F for rows 21-28 describing the
#   remaining data-structure changes needed for the 2-sex model.
# - Bold the new section headers ("Questions", "combined/separate
#   switches", "data structure changes").
# - Column A is widened to fit its (now longer) text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Clear out the old note that lived at E1:E4 -------------------------
$ws.Range("E1").Value = ""
$ws.Range("E2").Value = ""

# ---- Mark the survey/fishery composition rows (10-18) as done -----------
$ws.Range("B10").Value = "y"
$ws.Range("B11").Value = "y"
$ws.Range("B12").Value = "y"
$ws.Range("B13").Value = "y"
$ws.Range("B14").Value = "y"
$ws.Range("B15").Value = "y"
$ws.Range("B16").Value = "y"
$ws.Range("B17").Value = "y"
$ws.Range("B18").Value = "y"

# ---- Notes for rows 10-16 describing what was done -----------------------
# New strings get interned into the shared-string table in first-use order,
# so the order below matters to reproduce the exact sharedStrings layout.
$ws.Range("E11").Value = "wt"
$ws.Range("E12").Value = "empirical selectivity"
$ws.Range("E13").Value = "composition data"
$ws.Range("E10").Value = "combined/separate switches"
$ws.Range("E10").Font.Bold = $true

$ws.Range("E14").Value = "srv_comp_hat"
$ws.Range("E15").Value = "fsh_comp_hat"
$ws.Range("E16").Value = "UobsWtAge_hat"

# ---- New "data structure changes" question block (rows 21-28) ----------
$ws.Range("F22").Value = "sex in 4th column"
$ws.Range("E23").Value = "srv_biom"
$ws.Range("E24").Value = "fsh_biom"
$ws.Range("F23").Value = "remove sex column"
$ws.Range("E25").Value = "control"
$ws.Range("F25").Value = "nsex in 5th row"
$ws.Range("F26").Value = "R_sexr in 7th"
$ws.Range("F27").Value = "SSB_wt_index in 12"
$ws.Range("F28").Value = "spawn mo in 6"

$ws.Range("E21").Value = "data structure changes"
$ws.Range("E21").Font.Bold = $true

$ws.Range("E22").Value = "empirical selectivity"
$ws.Range("F24").Value = "remove sex column"

# ---- New "Questions" block (rows 3-4) ------------------------------------
$ws.Range("E3").Value = "Questions"
$ws.Range("E3").Font.Bold = $true
$ws.Range("E4").Value = "sex ratio in the likelihood?"

# ---- Column A autosize so the longer labels fit --------------------------
# (target best-fit width is 19.43 chars; engine quantizes ColumnWidth writes
# to 1/6-char steps, so 18.666... reliably lands on the nearest achievable
# column width of 19.5 chars)
$ws.Columns.Item(1).ColumnWidth = 18.666666666666668

# ---- Move the active selection -------------------------------------------
$ws.Range("L14").Select()
